$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

# Insert 3 new rows above row 26 (shifts old rows 26-28 down to 29-31,
# leaving rows 24-25 available for the new "power consumption by sector" block)
$ws.Rows("26:28").Insert()

# New row 24: ACT_BND / 2030 / 0 / DTCAR / TRAGSL
$ws.Range("D24").Value = "ACT_BND"
$ws.Range("E24").Value = 2030
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = "DTCAR"
$ws.Range("J24").Value = "TRAGSL"

# New row 25: ACT_BND / 0 / 2 / DTCAR / TRAGSL
$ws.Range("D25").Value = "ACT_BND"
$ws.Range("E25").Value = 0
$ws.Range("H25").Value = 2
$ws.Range("I25").Value = "DTCAR"
$ws.Range("J25").Value = "TRAGSL"

# Update selection / active cell to K24
$ws.Range("K24").Select()
